$wb = $excel.ActiveWorkbook

# Update "想去人数" (number of people interested) values on both the
# "展览" sheet and the "全部类型" sheet (which mirrors the same data).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 336
    $ws.Range("F4").Value = 72
}
